$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Relatorios"

# Row 2 "Diferença" value was negative (-R$ 90,00); correct it to positive (R$ 90,00)
$ws.Range("E2").Value = "R$ 90,00"

# Row 4 "Diferença" value was a flat repeat of "R$ 1.000,00"; it should be the
# negative difference (-R$ 1.000,00)
$ws.Range("E4").Value = "-R$ 1.000,00"

# Re-color the rows to reflect the new sign of their "Diferença" value:
# row 2 is now positive -> green (matches the color row 4 used to have)
# row 4 is now negative -> red (matches the color row 2 used to have)
$ws.Range("A2:E2").Interior.Color = 32768
$ws.Range("A4:E4").Interior.Color = 255

# Column E width nudges slightly wider (bestFit re-measure) to fit the new text
$ws.Columns.Item(5).ColumnWidth = 15.417
